$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped
# from 45203 (2023-10-04) to 45205 (2023-10-06) for every data row
# (rows 2 through 319).
$ws.Range("C2:C319").Value = 45205
